$wb = $excel.ActiveWorkbook

function Set-RowData($sheet, $rowNum, $rowHash) {
    foreach ($col in $rowHash.Keys) {
        $colIndex = switch ($col) {
            'B' { 2 }
            'C' { 3 }
            'D' { 4 }
            'E' { 5 }
        }
        $sheet.Cells.Item($rowNum, $colIndex).Value = $rowHash[$col]
    }
}

$sheet1Data = @{}
$sheet1Data[3] = @{'B'="Departmental DCA"; 'C'="Count"; 'D'="Costs"; 'E'="Proportion costs"}
$sheet1Data[4] = @{'B'="Green"; 'C'=2; 'D'=10273; 'E'=0.01831566485881267}
$sheet1Data[5] = @{'B'="Amber/Green"; 'C'=1; 'D'=67326; 'E'=0.1200350873439522}
$sheet1Data[6] = @{'B'="Amber"; 'C'=2; 'D'=498; 'E'=0.0008878809597672254}
$sheet1Data[7] = @{'B'="Amber/Red"; 'C'=1; 'D'=482789; 'E'=0.8607613668374678}
$sheet1Data[8] = @{'B'="Red"; 'C'=0; 'D'=0; 'E'=0}
$sheet1Data[9] = @{'B'="None"; 'C'=0; 'D'=0; 'E'=0}
$sheet1Data[10] = @{'B'="Total"; 'C'=6; 'D'=560886; 'E'=1}
$sheet1Data[13] = @{'B'="SRO Finance confidence"; 'C'="Count"; 'D'="Costs"; 'E'="Proportion costs"}
$sheet1Data[14] = @{'B'="Green"; 'C'=2; 'D'=10273; 'E'=0.01831566485881267}
$sheet1Data[15] = @{'B'="Amber/Green"; 'C'=0; 'D'=0; 'E'=0}
$sheet1Data[16] = @{'B'="Amber"; 'C'=4; 'D'=550613; 'E'=0.9816843351411874}
$sheet1Data[17] = @{'B'="Amber/Red"; 'C'=0; 'D'=0; 'E'=0}
$sheet1Data[18] = @{'B'="Red"; 'C'=0; 'D'=0; 'E'=0}
$sheet1Data[19] = @{'B'="None"; 'C'=0; 'D'=0; 'E'=0}
$sheet1Data[20] = @{'B'="Total"; 'C'=6; 'D'=560886; 'E'=1}
$sheet1Data[23] = @{'B'="SRO Benefits RAG"; 'C'="Count"; 'D'="Costs"; 'E'="Proportion costs"}
$sheet1Data[24] = @{'B'="Green"; 'C'=2; 'D'=4302; 'E'=0.007670007809073501}
$sheet1Data[25] = @{'B'="Amber/Green"; 'C'=0; 'D'=0; 'E'=0}
$sheet1Data[26] = @{'B'="Amber"; 'C'=4; 'D'=556584; 'E'=0.9923299921909265}
$sheet1Data[27] = @{'B'="Amber/Red"; 'C'=0; 'D'=0; 'E'=0}
$sheet1Data[28] = @{'B'="Red"; 'C'=0; 'D'=0; 'E'=0}
$sheet1Data[29] = @{'B'="None"; 'C'=0; 'D'=0; 'E'=0}
$sheet1Data[30] = @{'B'="Total"; 'C'=6; 'D'=560886; 'E'=1}
$sheet1Data[33] = @{'B'="SRO Schedule Confidence"; 'C'="Count"; 'D'="Costs"; 'E'="Proportion costs"}
$sheet1Data[34] = @{'B'="Green"; 'C'=2; 'D'=10273; 'E'=0.01831566485881267}
$sheet1Data[35] = @{'B'="Amber/Green"; 'C'=0; 'D'=0; 'E'=0}
$sheet1Data[36] = @{'B'="Amber"; 'C'=3; 'D'=550515; 'E'=0.9815096115788235}
$sheet1Data[37] = @{'B'="Amber/Red"; 'C'=0; 'D'=0; 'E'=0}
$sheet1Data[38] = @{'B'="Red"; 'C'=1; 'D'=98; 'E'=0.0001747235623638315}
$sheet1Data[39] = @{'B'="None"; 'C'=0; 'D'=0; 'E'=0}
$sheet1Data[40] = @{'B'="Total"; 'C'=6; 'D'=560886; 'E'=1}
$sheet1Data[43] = @{'B'="Overall Resource DCA - Now"; 'C'="Count"; 'D'="Costs"; 'E'="Proportion costs"}
$sheet1Data[44] = @{'B'="Green"; 'C'=0; 'D'=0; 'E'=0}
$sheet1Data[45] = @{'B'="Amber/Green"; 'C'=0; 'D'=0; 'E'=0}
$sheet1Data[46] = @{'B'="Amber"; 'C'=0; 'D'=0; 'E'=0}
$sheet1Data[47] = @{'B'="Amber/Red"; 'C'=0; 'D'=0; 'E'=0}
$sheet1Data[48] = @{'B'="Red"; 'C'=0; 'D'=0; 'E'=0}
$sheet1Data[49] = @{'B'="None"; 'C'=6; 'D'=560886; 'E'=1}
$sheet1Data[50] = @{'B'="Total"; 'C'=6; 'D'=560886; 'E'=1}

$sheet2Data = @{}
$sheet2Data[3] = @{'B'="Departmental DCA"; 'C'="Count"; 'D'="Costs"; 'E'="Proportion costs"}
$sheet2Data[4] = @{'B'="Green"; 'C'=0; 'D'=0; 'E'=0}
$sheet2Data[5] = @{'B'="Amber/Green"; 'C'=2; 'D'=1188; 'E'=0.1378030390905927}
$sheet2Data[6] = @{'B'="Amber"; 'C'=3; 'D'=7433; 'E'=0.8621969609094072}
$sheet2Data[7] = @{'B'="Amber/Red"; 'C'=0; 'D'=0; 'E'=0}
$sheet2Data[8] = @{'B'="Red"; 'C'=0; 'D'=0; 'E'=0}
$sheet2Data[9] = @{'B'="None"; 'C'=0; 'D'=0; 'E'=0}
$sheet2Data[10] = @{'B'="Total"; 'C'=5; 'D'=8621; 'E'=1}
$sheet2Data[13] = @{'B'="SRO Finance confidence"; 'C'="Count"; 'D'="Costs"; 'E'="Proportion costs"}
$sheet2Data[14] = @{'B'="Green"; 'C'=1; 'D'=672; 'E'=0.07794919382902216}
$sheet2Data[15] = @{'B'="Amber/Green"; 'C'=0; 'D'=0; 'E'=0}
$sheet2Data[16] = @{'B'="Amber"; 'C'=4; 'D'=7949; 'E'=0.9220508061709779}
$sheet2Data[17] = @{'B'="Amber/Red"; 'C'=0; 'D'=0; 'E'=0}
$sheet2Data[18] = @{'B'="Red"; 'C'=0; 'D'=0; 'E'=0}
$sheet2Data[19] = @{'B'="None"; 'C'=0; 'D'=0; 'E'=0}
$sheet2Data[20] = @{'B'="Total"; 'C'=5; 'D'=8621; 'E'=1}
$sheet2Data[23] = @{'B'="SRO Benefits RAG"; 'C'="Count"; 'D'="Costs"; 'E'="Proportion costs"}
$sheet2Data[24] = @{'B'="Green"; 'C'=0; 'D'=0; 'E'=0}
$sheet2Data[25] = @{'B'="Amber/Green"; 'C'=0; 'D'=0; 'E'=0}
$sheet2Data[26] = @{'B'="Amber"; 'C'=4; 'D'=8353; 'E'=0.9689131191277114}
$sheet2Data[27] = @{'B'="Amber/Red"; 'C'=0; 'D'=0; 'E'=0}
$sheet2Data[28] = @{'B'="Red"; 'C'=0; 'D'=0; 'E'=0}
$sheet2Data[29] = @{'B'="None"; 'C'=1; 'D'=268; 'E'=0.0310868808722886}
$sheet2Data[30] = @{'B'="Total"; 'C'=5; 'D'=8621; 'E'=1}
$sheet2Data[33] = @{'B'="Overall Resource DCA - Now"; 'C'="Count"; 'D'="Costs"; 'E'="Proportion costs"}
$sheet2Data[34] = @{'B'="Green"; 'C'=2; 'D'=1188; 'E'=0.1378030390905927}
$sheet2Data[35] = @{'B'="Amber/Green"; 'C'=0; 'D'=0; 'E'=0}
$sheet2Data[36] = @{'B'="Amber"; 'C'=3; 'D'=7433; 'E'=0.8621969609094072}
$sheet2Data[37] = @{'B'="Amber/Red"; 'C'=0; 'D'=0; 'E'=0}
$sheet2Data[38] = @{'B'="Red"; 'C'=0; 'D'=0; 'E'=0}
$sheet2Data[39] = @{'B'="None"; 'C'=0; 'D'=0; 'E'=0}
$sheet2Data[40] = @{'B'="Total"; 'C'=5; 'D'=8621; 'E'=1}

# --- Rename the existing (first) sheet and load it with the new Q4_19_20 data ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Q4_19_20"

foreach ($rowNum in $sheet1Data.Keys) {
    Set-RowData $ws1 $rowNum $sheet1Data[$rowNum]
}

# --- Add a new sheet right after it for the Q4_18_19 data (previous period refreshed) ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Q4_18_19"

foreach ($rowNum in $sheet2Data.Keys) {
    Set-RowData $ws2 $rowNum $sheet2Data[$rowNum]
}

# Keep the first sheet ("Q4_19_20") as the active tab
$ws1.Activate()
